$d = $word.ActiveDocument

$replacements = @(
    @("360÷6=", "661÷3="),
    @("907÷3=", "143÷8="),
    @("151÷7=", "774÷2="),
    @("727÷9=", "444÷4="),
    @("876÷8=", "310÷8="),
    @("445÷3=", "148÷7="),
    @("288÷8=", "579÷2="),
    @("595÷3=", "799÷5="),
    @("603÷6=", "294÷2="),
    @("882÷4=", "965÷6="),
    @("817÷3=", "993÷6="),
    @("147÷4=", "579÷4="),
    @("121÷6=", "142÷9="),
    @("308÷3=", "707÷4="),
    @("325÷5=", "707÷7="),
    @("106÷7=", "300÷6="),
    @("901÷8=", "121÷7="),
    @("222÷9=", "514÷4="),
    @("104÷2=", "953÷9="),
    @("774÷3=", "690÷3="),
    @("938÷7=", "848÷4="),
    @("711÷2=", "944÷4="),
    @("867÷9=", "189÷5="),
    @("975÷3=", "623÷2="),
    @("534÷8=", "831÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
